# Insert a new price-record row for "Artic Snow" / "Especial" at row 256.
# This pushes the existing rows 256-342 down to 257-343 (their content is
# preserved as-is), and the new row 256 is populated with the values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 256..342 down by one row, making room for the new record.
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256.
$ws.Range("A256").Value = 8
$ws.Range("B256").Value = "Terminal La Palmera de La Serena"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = 44636
$ws.Range("E256").Value = 4
$ws.Range("F256").Value = "Fruta"
$ws.Range("G256").Value = 100103
$ws.Range("H256").Value = "Frutos de hueso (carozo)"
$ws.Range("I256").Value = 100103006
$ws.Range("J256").Value = "Nectarín"
$ws.Range("K256").Value = "Artic Snow"
$ws.Range("L256").Value = "Especial"
$ws.Range("M256").Value = 20
$ws.Range("N256").Value = 390000
$ws.Range("O256").Value = 400000
$ws.Range("P256").Value = 395000
$ws.Range("Q256").Value = "`$/bins (420 kilos)"
$ws.Range("R256").Value = "Región de O'Higgins"
$ws.Range("S256").Value = 940
$ws.Range("T256").Value = 420
